$d = $word.ActiveDocument

# Locate the paragraph that ends with "...Give a person a warning" (the
# last paragraph of the "Cart features" note) — the diff appends five new
# paragraphs right after it, before the closing section properties.
$anchorRange = $d.Content
$found = $anchorRange.Find.Execute("Give a person a warning", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $anchorRange.Paragraphs(1)

# Paragraph 1: "Login features- use the same thing that we used in the RPS"
$anchorPara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "Login features- use the same thing that we used in the RPS"

# Paragraph 2: "Database – use the same thing we did in the RPS"
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "Database " + [char]0x2013 + " use the same thing we did in the RP"
$p2.Range.InsertAfter("S")

# Paragraph 3: the login/continue-shopping explanation
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "In RPS you "
$p3.Range.InsertAfter("can start a new game at the end, use that to ")
$p3.Range.InsertAfter("continue shopping ")
$p3.Range.InsertAfter("when you finish one.")
$p3.Range.InsertAfter(" If they don" + [char]0x2019 + "t want to play a new game they have to log out ")
$p3.Range.InsertAfter("and then")
$p3.Range.InsertAfter(" it goes back to login. So do the same ")
$p3.Range.InsertAfter("with the app, if they don" + [char]0x2019 + "t want to ")
$p3.Range.InsertAfter("continue shopping they can log out and then it goes to the log in prompt")

# Paragraph 4: the store-locations clarification
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Range.Text = "Clarification, you can have one store different "
$p4.Range.InsertAfter("locations or")
$p4.Range.InsertAfter(" each store is each location")
